$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'303.40"
$ws.Range("E2").Value = "'-0.90%"
$ws.Range("G2").Value = "'4"
$ws.Range("D3").Value = "'35.77"
$ws.Range("E3").Value = "'-1.05%"
$ws.Range("G3").Value = "'4"
$ws.Range("D4").Value = "'5.041"
$ws.Range("E4").Value = "'-1.27%"
$ws.Range("G4").Value = "'4"
$ws.Range("D5").Value = "'0.07896"
$ws.Range("E5").Value = "'-2.63%"
$ws.Range("G5").Value = "'4"
$ws.Range("D6").Value = "'1.838"
$ws.Range("E6").Value = "'-5.40%"
$ws.Range("G6").Value = "'4"
$ws.Range("D7").Value = "'4.103"
$ws.Range("E7").Value = "'-2.30%"
$ws.Range("G7").Value = "'4"
$ws.Range("D8").Value = "'7.782"
$ws.Range("E8").Value = "'0.14%"
$ws.Range("G8").Value = "'4"
$ws.Range("D9").Value = "'0.9198"
$ws.Range("E9").Value = "'-1.33%"
$ws.Range("G9").Value = "'4"
$ws.Range("D10").Value = "'0.1354"
$ws.Range("E10").Value = "'-3.71%"
$ws.Range("G10").Value = "'4"
$ws.Range("D11").Value = "'0.1884"
$ws.Range("E11").Value = "'-1.92%"
$ws.Range("G11").Value = "'4"
$ws.Range("D12").Value = "'0.09036"
$ws.Range("E12").Value = "'-1.89%"
$ws.Range("G12").Value = "'4"
$ws.Range("D13").Value = "'0.03472"
$ws.Range("E13").Value = "'-1.98%"
$ws.Range("G13").Value = "'4"
$ws.Range("D14").Value = "'0.09818"
$ws.Range("E14").Value = "'-0.28%"
$ws.Range("G14").Value = "'4"
$ws.Range("D15").Value = "'0.001408"
$ws.Range("E15").Value = "'-0.73%"
$ws.Range("G15").Value = "'4"
$ws.Range("D16").Value = "'0.006067"
$ws.Range("E16").Value = "'4.30%"
$ws.Range("G16").Value = "'4"
$ws.Range("D17").Value = "'3.724"
$ws.Range("E17").Value = "'3.55%"
$ws.Range("G17").Value = "'4"
$ws.Range("D18").Value = "'3.249"
$ws.Range("E18").Value = "'9.03%"
$ws.Range("G18").Value = "'4"
$ws.Range("D19").Value = "'0.3439"
$ws.Range("G19").Value = "'4"
$ws.Range("B20").Value = "ProBitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D20").Value = "'0.1343"
$ws.Range("E20").Value = "'-0.55%"
$ws.Range("G20").Value = "'4"
$ws.Range("B21").Value = "MCDex"
$ws.Range("C21").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D21").Value = "'5.193"
$ws.Range("E21").Value = "'6.16%"
$ws.Range("G21").Value = "'4"
$ws.Range("D22").Value = "'0.2195"
$ws.Range("E22").Value = "'-9.00%"
$ws.Range("G22").Value = "'4"
$ws.Range("D23").Value = "'0.04402"
$ws.Range("E23").Value = "'-2.39%"
$ws.Range("G23").Value = "'4"
$ws.Range("E24").Value = "'1.75%"
$ws.Range("G24").Value = "'4"
$ws.Range("D25").Value = "'0.004607"
$ws.Range("E25").Value = "'-5.48%"
$ws.Range("G25").Value = "'4"
$ws.Range("D26").Value = "'0.0001301"
$ws.Range("E26").Value = "'4.91%"
$ws.Range("G26").Value = "'4"
$ws.Range("D27").Value = "'0.0004450"
$ws.Range("E27").Value = "'0.19%"
$ws.Range("G27").Value = "'4"
$ws.Range("G28").Value = "'4"
$ws.Range("G29").Value = "'4"
$ws.Range("G30").Value = "'4"
$ws.Range("G31").Value = "'4"
$ws.Range("G32").Value = "'4"
$ws.Range("G33").Value = "'4"
$ws.Range("G34").Value = "'4"
$ws.Range("G35").Value = "'4"
$ws.Range("G36").Value = "'4"
$ws.Range("G37").Value = "'4"
$ws.Range("G38").Value = "'4"
$ws.Range("D39").Value = "'0.01929"
$ws.Range("E39").Value = "'-3.92%"
$ws.Range("G39").Value = "'4"
$ws.Range("D40").Value = "'0.05070"
$ws.Range("G40").Value = "'4"
$ws.Range("D41").Value = "'0.007616"
$ws.Range("E41").Value = "'-0.56%"
$ws.Range("G41").Value = "'4"
$ws.Range("D42").Value = "'0.01014"
$ws.Range("E42").Value = "'-8.03%"
$ws.Range("G42").Value = "'4"
$ws.Range("D43").Value = "'0.1338"
$ws.Range("E43").Value = "'-3.09%"
$ws.Range("G43").Value = "'4"
$ws.Range("D44").Value = "'0.002152"
$ws.Range("E44").Value = "'2.45%"
$ws.Range("G44").Value = "'4"
$ws.Range("D45").Value = "'0.01020"
$ws.Range("E45").Value = "'-3.64%"
$ws.Range("G45").Value = "'4"
$ws.Range("D46").Value = "'0.00006154"
$ws.Range("E46").Value = "'-4.78%"
$ws.Range("G46").Value = "'4"
$ws.Range("D47").Value = "'0.00000000752"
$ws.Range("E47").Value = "'0.03%"
$ws.Range("G47").Value = "'4"
$ws.Range("G48").Value = "'4"
$ws.Range("D49").Value = "'0.001662"
$ws.Range("E49").Value = "'39.39%"
$ws.Range("G49").Value = "'4"
$ws.Range("D50").Value = "'0.00002104"
$ws.Range("E50").Value = "'0.03%"
$ws.Range("G50").Value = "'4"
$ws.Range("D51").Value = "'0.0002004"
$ws.Range("E51").Value = "'0.03%"
$ws.Range("G51").Value = "'4"
